# Add tests for nested select one, select multiple and gps data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Insert 3 new rows after the "hiv medication" question (before "begin repeat") ---
$ws.Range("A7:A9").EntireRow.Insert()

$ws.Range("A7").Value = "geopoint"
$ws.Range("B7").Value = "gps"
$ws.Range("C7").Value = "Input gps coordinates"
$ws.Rows.Item(7).RowHeight = 24

$ws.Range("A8").Value = "select_multiple pizza_toppings"
$ws.Range("B8").Value = "pizza_topping"
$ws.Range("C8").Value = "What toppings to do you prefer?"
$ws.Rows.Item(8).RowHeight = 24

$ws.Range("A9").Value = "select_one yes_no"
$ws.Range("B9").Value = "customer_satisfaction"
$ws.Range("C9").Value = "Are you satisfied with the level of service received?"
$ws.Rows.Item(9).RowHeight = 24

# --- Insert 2 new rows inside the repeat group, right after "Last name" (before "integer") ---
$ws.Range("A14:A15").EntireRow.Insert()

$ws.Range("A14").Value = "geopoint"
$ws.Range("B14").Value = "gps"
$ws.Range("C14").Value = "Input gps coordinates"
$ws.Rows.Item(14).RowHeight = 24

$ws.Range("A15").Value = "select_multiple pizza_toppings"
$ws.Range("B15").Value = "pizza_topping"
$ws.Range("C15").Value = "What toppings to do you prefer?"
$ws.Rows.Item(15).RowHeight = 24

# --- Insert 1 new row inside the repeat group, right after "age" (before "end group") ---
$ws.Range("A17:C17").EntireRow.Insert()

$ws.Range("A17").Value = "select_one yes_no"
$ws.Range("B17").Value = "customer_satisfaction"
$ws.Range("C17").Value = "Are you satisfied with the level of service received?"
$ws.Rows.Item(17).RowHeight = 24

# --- Column width tweaks on the survey sheet ---
$ws.Columns.Item(1).ColumnWidth = 37.92
$ws.Columns.Item(2).ColumnWidth = 25.84

# --- Selection / active sheet bookkeeping ---
$ws.Range("A15").Select()
$ws.Activate()
